$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07).
$ws.Range("C2:C476").Value = 45206
